$d = $word.ActiveDocument

# --- Helper: plain whole-phrase replace (used where run formatting is
# uniform across the matched text, so Find/Replace losing run splits is
# harmless). wdFindContinue = 1, wdReplaceOne = 1 (we pass explicit args
# matching the signature shown in the task prompt).
function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- Helper: turn "<n>st" into "<n>0th" in place while preserving the
# superscript run's character formatting (only the visible digits/letters
# change; a new plain "0" run is inserted ahead of the still-superscript
# "th"). $context is a unique surrounding phrase that still contains the
# literal "st" to be fixed; $stOffset is the character offset of that
# "st" within $context.
function Fix-Ordinal($context, $stOffset) {
    $hit = $d.Content.Duplicate
    $found = $hit.Find.Execute($context, $false, $false, $false, $false, `
                                $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    $stStart = $hit.Start + $stOffset
    $stRange = $d.Range($stStart, $stStart + 2)
    if ($stRange.Text -ne "st") {
        return
    }
    $stRange.InsertBefore("0")
    $thRange = $d.Range($stStart + 1, $stStart + 3)
    $thRange.Text = "th"
}

# Date: 1st February 2024  ->  Date: 10th February 2024
Fix-Ordinal "Date: 1st February 2024" 7

# Your employment will begin on 1st February 2024 -> ... 10th February 2024
Fix-Ordinal "begin on 1st February 2024" 10

# Employee name: Sagar Chaudhari -> Ganesh Mali (Name field, Dear field,
# and the signature line all share this exact phrase).
Replace-Text "Sagar Chaudhari" "Ganesh Mali"

# Address: H. No: 580/14, Cortalim , Goa -> Nerul , Goa
Replace-Text "H. No: 580/14, Cortalim , Goa" "Nerul , Goa"

# Role: Technical Assistant (Fresher) -> Driver cum Helper (appears in
# the "Appointed as" field and in the offer paragraph).
Replace-Text "Technical Assistant (Fresher)" "Driver cum Helper"

# Monthly CTC: 9000 -> 15000
Replace-Text "9000" "15000"

Write-Output "edits applied"
